$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row (A:R only) before row 118, shifting the existing row 118 data down to row 119
$ws.Range("A118:R118").Insert(-4121)  # xlShiftDown

# Copy formatting for the new row 118 from row 117 (keeps date style etc. consistent)
$ws.Range("A117:R117").Copy()
$ws.Range("A118:R118").PasteSpecial(-4122)  # xlPasteFormats

# Fill new row 118 with the new week's data
$ws.Range("A118").Value = 8
$ws.Range("B118").Value = "Terminal La Palmera de La Serena"
$ws.Range("C118").Value = "Coquimbo"
$ws.Range("D118").Value = 44595
$ws.Range("E118").Value = 4
$ws.Range("F118").Value = 100112044
$ws.Range("G118").Value = "Perejil"
$ws.Range("H118").Value = "Sin especificar"
$ws.Range("I118").Value = "Primera"
$ws.Range("J118").Value = 3000
$ws.Range("K118").Value = 2500
$ws.Range("L118").Value = 2800
$ws.Range("M118").Value = 2650
$ws.Range("N118").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O118").Value = "Provincia del Elquí"
$ws.Range("P118").Value = 1767
$ws.Range("Q118").Value = 1.5
$ws.Range("R118").Value = "Hortaliza"
